$wb = $excel.ActiveWorkbook

# Rename the "inflow" sheet to "inflow1".
$wsInflow1 = $wb.Worksheets.Item("inflow")
$wsInflow1.Name = "inflow1"

$wsInflow2 = $wb.Worksheets.Item("inflow2")
$wsOutflow = $wb.Worksheets.Item("outflow")
$wsBypass  = $wb.Worksheets.Item("bypass")

# The shared string "flow_rate_value" is retired in favor of "flow" -- update the
# header cells on the sheets that still used the old label.
$wsOutflow.Range("B1").Value = "flow"
$wsBypass.Range("B1").Value = "flow"

# Restore each sheet's remembered cursor/selection position.
$wsInflow1.Activate()
$null = $wsInflow1.Range("D25").Select()

$wsInflow2.Activate()
$null = $wsInflow2.Range("D38").Select()

$wsOutflow.Activate()
$null = $wsOutflow.Range("B2").Select()

# "bypass" ends up as the active/selected sheet.
$wsBypass.Activate()
$null = $wsBypass.Range("A2").Select()
